$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.05
$ws.Range("I2").Value = 3.9
$ws.Range("N6").Value = 2.04
$ws.Range("O6").Value = 1.86
$ws.Range("N7").Value = 2.05
$ws.Range("O7").Value = 1.68
$ws.Range("N8").Value = 1.62
$ws.Range("O8").Value = 2.15
$ws.Range("L9").Value = 1.36
$ws.Range("M9").Value = 3.2
$ws.Range("N9").Value = 2.05
$ws.Range("O9").Value = 1.68
$ws.Range("N13").Value = 1.73
$ws.Range("O13").Value = 2.08
$ws.Range("N15").Value = 1.5
$ws.Range("N16").Value = 1.5
$ws.Range("N17").Value = 2.3
$ws.Range("O17").Value = 1.62
$ws.Range("J19").Value = 1.05
$ws.Range("K19").Value = 11
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1.9
$ws.Range("K21").Value = 10
$ws.Range("J23").Value = 1.13
$ws.Range("K23").Value = 6
$ws.Range("AH23").Value = 67
$ws.Range("N50").Value = 2.03
$ws.Range("O50").Value = 1.83
$ws.Range("J56").Value = 1.03
$ws.Range("K56").Value = 15
$ws.Range("L56").Value = 1.2
$ws.Range("M56").Value = 4.33
$ws.Range("N56").Value = 1.67
$ws.Range("O56").Value = 2.15
$ws.Range("L57").Value = 1.17
$ws.Range("M57").Value = 5
$ws.Range("N57").Value = 1.6
$ws.Range("O57").Value = 2.3
$ws.Range("J59").Value = 1.07
$ws.Range("K59").Value = 9
$ws.Range("P64").Value = 1.33
$ws.Range("Q64").Value = 3.25
$ws.Range("R64").Value = 1.57
$ws.Range("S64").Value = 2.25
$ws.Range("T64").Value = 11
$ws.Range("U64").Value = 15
$ws.Range("Y64").Value = 23
$ws.Range("AB64").Value = 12
$ws.Range("AD64").Value = 126
$ws.Range("AE64").Value = 11
$ws.Range("AF64").Value = 15
$ws.Range("AJ64").Value = 23
$ws.Range("N71").Value = 2.05
$ws.Range("G77").Value = 2.35
$ws.Range("I77").Value = 2.88
$ws.Range("L77").Value = 1.33
$ws.Range("M77").Value = 3.25
$ws.Range("P77").Value = 1.44
$ws.Range("Q77").Value = 2.63
$ws.Range("U77").Value = 11
$ws.Range("W77").Value = 21
$ws.Range("X77").Value = 19
$ws.Range("AF77").Value = 15
$ws.Range("AJ77").Value = 34
$ws.Range("G79").Value = 4.1
$ws.Range("P79").Value = 1.33
$ws.Range("Q79").Value = 3.25
$ws.Range("R79").Value = 1.67
$ws.Range("S79").Value = 2.1
$ws.Range("U79").Value = 21
$ws.Range("V79").Value = 13
$ws.Range("X79").Value = 29
$ws.Range("Z79").Value = 13
$ws.Range("AB79").Value = 13
$ws.Range("AD79").Value = 151
$ws.Range("AE79").Value = 8.5
$ws.Range("AF79").Value = 9.5
$ws.Range("AI79").Value = 13
$ws.Range("G84").Value = 2.63
$ws.Range("I84").Value = 2.4
$ws.Range("V84").Value = 10
$ws.Range("W84").Value = 26
$ws.Range("Y84").Value = 26
$ws.Range("AE84").Value = 9.5
$ws.Range("AF84").Value = 13
$ws.Range("AH84").Value = 23
$ws.Range("AI84").Value = 19
$ws.Range("AJ84").Value = 26
$ws.Range("K85").Value = 8
$ws.Range("N85").Value = 2.4
$ws.Range("O85").Value = 1.53
$ws.Range("N88").Value = 2.25
$ws.Range("O88").Value = 1.62
$ws.Range("G89").Value = 1.3
$ws.Range("H89").Value = 6
$ws.Range("I89").Value = 8
$ws.Range("J89").Value = 1.01
$ws.Range("K89").Value = 26
$ws.Range("N89").Value = 1.33
$ws.Range("O89").Value = 3.4
$ws.Range("R89").Value = 1.62
$ws.Range("S89").Value = 2.2
$ws.Range("U89").Value = 9
$ws.Range("W89").Value = 10
$ws.Range("Y89").Value = 19
$ws.Range("Z89").Value = 26
$ws.Range("AB89").Value = 17
$ws.Range("AF89").Value = 41
$ws.Range("AG89").Value = 23
$ws.Range("AH89").Value = 81
$ws.Range("J91").Value = 1.05
$ws.Range("K91").Value = 11
$ws.Range("N91").Value = 1.98
$ws.Range("O91").Value = 1.88
$ws.Range("J93").Value = 1.02
$ws.Range("L93").Value = 1.13
$ws.Range("J94").Value = 1.02
$ws.Range("L94").Value = 1.17
$ws.Range("N98").Value = 1.9
$ws.Range("O98").Value = 1.95
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 3.3
$ws.Range("I106").Value = 1.75
$ws.Range("K106").Value = 10
$ws.Range("L106").Value = 1.29
$ws.Range("M106").Value = 3.5
$ws.Range("N106").Value = 2
$ws.Range("O106").Value = 1.85
$ws.Range("T106").Value = 12
$ws.Range("Z106").Value = 10
$ws.Range("AA106").Value = 7
$ws.Range("AG106").Value = 8.5
$ws.Range("G127").Value = 1.9
$ws.Range("I127").Value = 4
$ws.Range("W127").Value = 15
$ws.Range("Z127").Value = 8.5
$ws.Range("AE127").Value = 10
$ws.Range("J128").Value = 1.06
$ws.Range("K128").Value = 10
$ws.Range("V128").Value = 12
$ws.Range("AH128").Value = 19
$ws.Range("G129").Value = 2.38
$ws.Range("I129").Value = 2.8
$ws.Range("N129").Value = 2.05
$ws.Range("O129").Value = 1.75
$ws.Range("AF129").Value = 15
$ws.Range("N130").Value = 2.08
$ws.Range("O130").Value = 1.73
$ws.Range("G133").Value = 2.63
$ws.Range("H133").Value = 2.9
$ws.Range("U133").Value = 12
$ws.Range("AA133").Value = 5.5
$ws.Range("AC133").Value = 51
$ws.Range("AG133").Value = 12
$ws.Range("H136").Value = 3.75
$ws.Range("K136").Value = 12
$ws.Range("P136").Value = 1.33
$ws.Range("Q136").Value = 3.25
$ws.Range("R136").Value = 1.7
$ws.Range("S136").Value = 2.05
$ws.Range("AJ136").Value = 34
$ws.Range("J149").Value = 1.08
$ws.Range("K149").Value = 8
$ws.Range("N151").Value = 1.93
$ws.Range("O151").Value = 1.93
$ws.Range("N152").Value = 1.83
$ws.Range("O152").Value = 2.03
